$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.906.36"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.545.22"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'205.48"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'0.484"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'21.30"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "1.765.49"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "1.555.13"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "26.886.24"
$ws.Range("D17").Value = "'61.57"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "0.0₃0682"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'7.17"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'4.02"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").Value = "'153.11"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'0.0457"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "1.361.73"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.973"
$ws.Range("E36").Value = "  +5.97%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'0.988"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.47"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "'2.22"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "'63.31"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'1.74"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "1.679.04"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").Value = "'86.22"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "'0.0506"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "0.0₇0968"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.0949"
$ws.Range("E51").Value = "  -0.24%  "
